# Rename the worksheet: StringLocalizations_BasicText -> Localization
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Localization"

# Append the two new localization keys as new rows (A134, A135) in column A only,
# mirroring the existing "key"-only rows already present at the bottom of the sheet.
$ws.Cells.Item(134, 1).Value = "BRANDING_FORCE_NAME"
$ws.Cells.Item(135, 1).Value = "BRANDING_FORCE_STRAPLINE"

# Move the view down to the newly added rows and select the last one, matching
# the author's final cursor position/scroll state when they made the edit.
$ws.Range("A135").Select()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 119
